$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H138").Value = 4088.4
$ws.Range("I138").Value = 1727.5294
$ws.Range("J138").Value = 4571.9517
$ws.Range("K138").Value = 5182.5882
$ws.Range("L138").Value = 13715.8551
$ws.Range("M138").Value = -42.58820000000014
$ws.Range("N138").Value = -23995.8551
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 5052.5
$ws.Range("I3").Value = 105
$ws.Range("K3").Value = 105
$ws.Range("M3").Value = 10
$ws.Range("H5").Value = 230
$ws.Range("I5").Value = 195
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 195
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -83
$ws.Range("N5").Value = -524
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H32").Value = 15310.87
$ws.Range("I32").Value = 12233.835
$ws.Range("J32").Value = 26886.38
$ws.Range("K32").Value = 12233.835
$ws.Range("L32").Value = 26886.38
$ws.Range("M32").Value = -11946.835
$ws.Range("N32").Value = -27460.38
$ws.Range("H45").Value = 2088.2727
$ws.Range("I45").Value = 2464.4
$ws.Range("J45").Value = 1774.8334
$ws.Range("K45").Value = 2464.4
$ws.Range("L45").Value = 1774.8334
$ws.Range("M45").Value = -2087.4
$ws.Range("N45").Value = -2528.8334
$ws.Range("H132").Value = 25302.63
$ws.Range("I132").Value = 32188.97
$ws.Range("K132").Value = 96566.91
$ws.Range("M132").Value = -94036.91
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 230
$ws.Range("I4").Value = 195
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 195
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -80
$ws.Range("N4").Value = -530
$ws.Range("H8").Value = 1082.25
$ws.Range("I8").Value = 1082.25
$ws.Range("K8").Value = 1082.25
$ws.Range("M8").Value = -942.25
$ws.Range("H81").Value = 36652.5
$ws.Range("J81").Value = 36652.5
$ws.Range("L81").Value = 36652.5
$ws.Range("N81").Value = -38774.5
$ws.Range("H82").Value = 13383.048
$ws.Range("J82").Value = 29287.5
$ws.Range("L82").Value = 29287.5
$ws.Range("N82").Value = -30053.5
$ws.Range("H84").Value = 36652.5
$ws.Range("J84").Value = 36652.5
$ws.Range("L84").Value = 109957.5
$ws.Range("N84").Value = -120565.5
$ws.Range("H85").Value = 13383.048
$ws.Range("J85").Value = 29287.5
$ws.Range("L85").Value = 29287.5
$ws.Range("N85").Value = -31939.5
$ws.Range("H107").Value = 1636.3948
$ws.Range("I107").Value = 1502.6086
$ws.Range("J107").Value = 1841.5333
$ws.Range("K107").Value = 1502.6086
$ws.Range("L107").Value = 1841.5333
$ws.Range("M107").Value = 417.3914
$ws.Range("N107").Value = -5681.5333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2540.233
$ws.Range("I31").Value = 1288.9117
$ws.Range("J31").Value = 3631.1282
$ws.Range("K31").Value = 1288.9117
$ws.Range("L31").Value = 3631.1282
$ws.Range("M31").Value = -993.9117000000001
$ws.Range("N31").Value = -4221.1282
$ws.Range("H34").Value = 2540.233
$ws.Range("I34").Value = 1288.9117
$ws.Range("J34").Value = 3631.1282
$ws.Range("K34").Value = 1288.9117
$ws.Range("L34").Value = 3631.1282
$ws.Range("M34").Value = -1086.9117
$ws.Range("N34").Value = -4035.1282
$ws.Range("H86").Value = 230482
$ws.Range("I86").Value = 314731.5
$ws.Range("K86").Value = 314731.5
$ws.Range("M86").Value = -313608.5
$ws.Range("H89").Value = 230482
$ws.Range("I89").Value = 314731.5
$ws.Range("K89").Value = 1573657.5
$ws.Range("M89").Value = -1568041.5
$ws.Range("H107").Value = 1254.1111
$ws.Range("I107").Value = 1274.9
$ws.Range("K107").Value = 1274.9
$ws.Range("M107").Value = 645.0999999999999
$ws.Range("H122").Value = 899.6842
$ws.Range("I122").Value = 739.2857
$ws.Range("J122").Value = 993.25
$ws.Range("K122").Value = 2217.8571
$ws.Range("L122").Value = 2979.75
$ws.Range("M122").Value = 232.1428999999998
$ws.Range("N122").Value = -7879.75
$ws.Range("H132").Value = 2321.0476
$ws.Range("I132").Value = 1775.4546
$ws.Range("J132").Value = 2921.2
$ws.Range("K132").Value = 5326.3638
$ws.Range("L132").Value = 8763.599999999999
$ws.Range("M132").Value = -2796.3638
$ws.Range("N132").Value = -13823.6
$ws.Range("H134").Value = 3039.3823
$ws.Range("I134").Value = 3260.5173
$ws.Range("J134").Value = 1756.8
$ws.Range("K134").Value = 9781.5519
$ws.Range("L134").Value = 5270.4
$ws.Range("M134").Value = -7246.5519
$ws.Range("N134").Value = -10340.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 685.6984
$ws.Range("I5").Value = 517.0213
$ws.Range("J5").Value = 1181.1875
$ws.Range("K5").Value = 1551.0639
$ws.Range("L5").Value = 3543.5625
$ws.Range("M5").Value = -1439.0639
$ws.Range("N5").Value = -3767.5625
$ws.Range("H122").Value = 1010.9697
$ws.Range("I122").Value = 406
$ws.Range("J122").Value = 1404.2
$ws.Range("K122").Value = 3654
$ws.Range("L122").Value = 12637.8
$ws.Range("M122").Value = -1204
$ws.Range("N122").Value = -17537.8
$ws.Range("H132").Value = 827.5526
$ws.Range("I132").Value = 401.88
$ws.Range("J132").Value = 1646.1538
$ws.Range("K132").Value = 3616.92
$ws.Range("L132").Value = 14815.3842
$ws.Range("M132").Value = -1086.92
$ws.Range("N132").Value = -19875.3842
$ws.Range("H135").Value = 685.6984
$ws.Range("I135").Value = 517.0213
$ws.Range("J135").Value = 1181.1875
$ws.Range("K135").Value = 4653.1917
$ws.Range("L135").Value = 10630.6875
$ws.Range("M135").Value = -2118.1917
$ws.Range("N135").Value = -15700.6875
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 31500000
$ws.Range("I3").Value = 31500000
$ws.Range("K3").Value = 31500000
$ws.Range("M3").Value = -31499884
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2565000.8
$ws.Range("I22").Value = 4762469
$ws.Range("J22").Value = 1288.1666
$ws.Range("K22").Value = 4762469
$ws.Range("L22").Value = 1288.1666
$ws.Range("M22").Value = -4762174
$ws.Range("N22").Value = -1878.1666
$ws.Range("H27").Value = 2565000.8
$ws.Range("I27").Value = 4762469
$ws.Range("J27").Value = 1288.1666
$ws.Range("K27").Value = 4762469
$ws.Range("L27").Value = 1288.1666
$ws.Range("M27").Value = -4762362
$ws.Range("N27").Value = -1502.1666
$ws.Range("H137").Value = 45843
$ws.Range("I137").Value = 40429
$ws.Range("J137").Value = 46925.8
$ws.Range("K137").Value = 40429
$ws.Range("L137").Value = 46925.8
$ws.Range("M137").Value = -35329
$ws.Range("N137").Value = -57125.8
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 1292.5
$ws.Range("I58").Value = 1292.5
$ws.Range("K58").Value = 1292.5
$ws.Range("M58").Value = -984.5
$ws.Range("H75").Value = 25000
$ws.Range("J75").Value = 25000
$ws.Range("L75").Value = 25000
$ws.Range("N75").Value = -26872
$ws.Range("H78").Value = 25000
$ws.Range("J78").Value = 25000
$ws.Range("L78").Value = 75000
$ws.Range("N78").Value = -84360
